$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at row 23, shifting existing rows 23-36 down to 24-37
# (same as the Nomura broker being inserted just under the existing
# "INST-FI" row rather than appended at the end).
$ws.Rows.Item(23).Insert()

# Populate the new row 23 with the Nomura broker record (NOMX-FI) — same
# counterparty details as INST-FI (row 22), just filed under its own
# broker id, plus a note in column G explaining the relationship.
$ws.Cells.Item(23, 1).Value = "NOMX-FI"
$ws.Cells.Item(23, 2).Value = "NOMURA INTERNATIONAL (HK) LIMITED"
$ws.Cells.Item(23, 3).Value = "NOMAGB2LXXX"
$ws.Cells.Item(23, 4).Value = "EUROCLEAR"

# "90997" looks numeric, but must be stored as text (shared string) like
# the rest of the Participant ID column — force text format before the
# write, then drop the number-format override so the cell keeps the
# workbook's default (unstyled) look, matching its neighbours.
$ws.Cells.Item(23, 5).NumberFormat = "@"
$ws.Cells.Item(23, 5).Value = "90997"
$ws.Cells.Item(23, 5).ClearFormats()

$ws.Cells.Item(23, 7).Value = "NOMX-FI is INST-FI"

# Match the author's final selection/cursor position.
$ws.Range("I23").Select()
